$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param([int]$r, [string]$timestamp, [string]$sender, [double]$senderId, [string]$phone, [string]$message, [string]$media)

    $ws.Range("A$r").Value = $timestamp
    $ws.Range("B$r").Value = $sender
    $ws.Range("C$r").Value = $senderId

    $ws.Range("D$r").Value = "'" + $phone
    $ws.Range("D$r").Style = "Normal"

    $ws.Range("E$r").Value = $message

    if ($media -ne "") {
        $ws.Range("F$r").Value = $media
    } else {
        $ws.Range("F$r").Style = "Normal"
    }

    $ws.Range("G$r").Style = "Normal"
}

Set-Row 33 "2025-10-02 14:11:48" "Noah" 8450689526 "13052054965" "Hey man what’s up?" ""
Set-Row 34 "2025-10-02 14:41:25" "Noah" 8450689526 "13052054965" "Hi" ""
Set-Row 35 "2025-10-02 14:43:50" "Noah" 8450689526 "13052054965" "Hey man" ""
Set-Row 36 "2025-10-02 14:44:07" "Noah" 8450689526 "13052054965" "What’s up man?" ""
Set-Row 37 "2025-10-02 14:45:00" "Noah" 8450689526 "13052054965" "Hey man" ""
Set-Row 38 "2025-10-02 14:50:33" "Noah" 8450689526 "13052054965" "Test messages" ""
Set-Row 39 "2025-10-02 14:51:56" "Noah" 8450689526 "13052054965" "Yo" ""
Set-Row 40 "2025-10-02 15:39:53" "Noah" 8450689526 "13052054965" "Hey man" ""
Set-Row 41 "2025-10-02 15:41:27" "Noah" 8450689526 "13052054965" "Test messages" ""
Set-Row 42 "2025-10-02 15:42:18" "Noah" 8450689526 "13052054965" "This is a test message" "my-node-server/public/uploads/images\photo_2025-10-02_19-42-18.jpg"
